$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-16 Monday" "2024-09-17 Tuesday"

Replace-Text "299×7=2093" "195×7=1365"
Replace-Text "531×4=2124" "552×8=4416"
Replace-Text "765×7=5355" "709×9=6381"
Replace-Text "759×9=6831" "248×6=1488"
Replace-Text "531×6=3186" "234×7=1638"

Replace-Text "459×6=2754" "305×9=2745"
Replace-Text "683×3=2049" "260×4=1040"
Replace-Text "687×4=2748" "409×3=1227"
Replace-Text "826×5=4130" "314×7=2198"
Replace-Text "978×2=1956" "149×4=596"

Replace-Text "827×3=2481" "243×5=1215"
Replace-Text "933×7=6531" "386×2=772"
Replace-Text "402×7=2814" "182×7=1274"
Replace-Text "357×7=2499" "976×5=4880"
Replace-Text "285×4=1140" "778×5=3890"

Replace-Text "884×4=3536" "656×2=1312"
Replace-Text "693×3=2079" "554×3=1662"
Replace-Text "946×8=7568" "120×4=480"
Replace-Text "344×5=1720" "589×7=4123"
Replace-Text "560×4=2240" "930×4=3720"

Replace-Text "601×8=4808" "856×2=1712"
Replace-Text "901×3=2703" "457×4=1828"
Replace-Text "494×7=3458" "514×8=4112"
Replace-Text "921×4=3684" "613×6=3678"
Replace-Text "150×5=750" "187×4=748"
